$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Weekly data refresh: rotate Fecha/Volumen/Precio values among existing rows
# Row 2 <- (old row 7 values)
$ws.Range("D2").Value = 45141
$ws.Range("J2").Value = 50
$ws.Range("K2").Value = 8500
$ws.Range("L2").Value = 9000
$ws.Range("M2").Value = 8800
$ws.Range("P2").Value = 587

# Row 3 <- (old row 8 values)
$ws.Range("D3").Value = 45119
$ws.Range("J3").Value = 50
$ws.Range("K3").Value = 20000
$ws.Range("L3").Value = 20000
$ws.Range("M3").Value = 20000
$ws.Range("P3").Value = 1333

# Row 5 <- (old row 2 values)
$ws.Range("D5").Value = 44749
$ws.Range("J5").Value = 90
$ws.Range("K5").Value = 17000
$ws.Range("L5").Value = 18000
$ws.Range("M5").Value = 17556
$ws.Range("P5").Value = 1170

# Row 7 <- (old row 5 values)
$ws.Range("D7").Value = 45091
$ws.Range("J7").Value = 40
$ws.Range("K7").Value = 20000
$ws.Range("L7").Value = 22000
$ws.Range("M7").Value = 21000
$ws.Range("P7").Value = 1400

# Row 8 <- (old row 10 values)
$ws.Range("D8").Value = 45084
$ws.Range("J8").Value = 90
$ws.Range("K8").Value = 22000
$ws.Range("L8").Value = 23000
$ws.Range("M8").Value = 22556
$ws.Range("P8").Value = 1504

# Row 10 <- (old row 3 values)
$ws.Range("D10").Value = 44839
$ws.Range("J10").Value = 50
$ws.Range("K10").Value = 15000
$ws.Range("L10").Value = 16000
$ws.Range("M10").Value = 15600
$ws.Range("P10").Value = 1040
